# Tweak to TNM values, running Balance model with different BaseMigs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CMAP Region, 2050 NetMigration: 375000 -> 350000
$ws.Range("C7").Value = 350000

# External IN, 2040 NetMigration: 5000 -> 8000
$ws.Range("C23").Value = 8000

# External IN, 2045 NetMigration: 0 -> 6000
$ws.Range("C24").Value = 6000

# External IN, 2050 NetMigration: 0 -> 5000
$ws.Range("C25").Value = 5000

# Update the active selection to match the author's final cursor position
$ws.Range("E21").Select()
